$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.247.72"
$ws.Range("E2").Value = "  +1.08%  "

$ws.Range("D3").Value = "1.798.64"
$ws.Range("E3").Value = "  +2.17%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "338.71"
$ws.Range("E5").Value = "  +0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9990"
$ws.Range("E6").Value = "  -0.13%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4782"
$ws.Range("E7").Value = "  +27.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3644"
$ws.Range("E8").Value = "  +8.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.50"
$ws.Range("E9").Value = "  -0.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07695"
$ws.Range("E10").Value = "  +7.66%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.144"
$ws.Range("E11").Value = "  +2.52%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.58"
$ws.Range("E12").Value = "  +1.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.001"
$ws.Range("E13").Value = "  -0.05%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.295"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.301"
$ws.Range("E15").Value = "  +2.30%  "

$ws.Range("D16").Value = "1.795.00"
$ws.Range("E16").Value = "  +2.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001094"
$ws.Range("E17").Value = "  +4.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06713"
$ws.Range("E18").Value = "  +2.21%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.72"
$ws.Range("E19").Value = "  +2.17%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9986"
$ws.Range("E20").Value = "  -0.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.30"
$ws.Range("E21").Value = "  +2.53%  "

$ws.Range("E22").Value = "  +2.51%  "

$ws.Range("D23").Value = "28.242.29"
$ws.Range("E23").Value = "  +1.09%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.01"
$ws.Range("E24").Value = "  +3.04%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.399"
$ws.Range("E25").Value = "  +0.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.57"
$ws.Range("E26").Value = "  +4.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.404"
$ws.Range("E27").Value = "  +3.62%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.13"
$ws.Range("E28").Value = "  -0.63%  "

$ws.Range("D29").Value = "2.002.93"
$ws.Range("E29").Value = "  +2.30%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.95"
$ws.Range("E30").Value = "  +2.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.270"
$ws.Range("E31").Value = "  +1.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.043"
$ws.Range("E32").Value = "  +0.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09650"
$ws.Range("E33").Value = "  +10.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.920"
$ws.Range("E34").Value = "  +2.67%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02383"
$ws.Range("E35").Value = "  +2.32%  "

$ws.Range("E36").Value = "  +0.03%  "

$ws.Range("E37").Value = "  +2.08%  "

$ws.Range("E38").Value = "  +2.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.208"
$ws.Range("E39").Value = "  +1.63%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2172"
$ws.Range("E40").Value = "  +3.43%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.482"
$ws.Range("E41").Value = "  +2.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.216"
$ws.Range("E42").Value = "  +0.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.084"
$ws.Range("E43").Value = "  +0.77%  "

$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9987"
$ws.Range("E44").Value = "  -0.17%  "

$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.08"
$ws.Range("E45").Value = "  +3.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.865"
$ws.Range("E46").Value = "  +0.95%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6125"
$ws.Range("E47").Value = "  +2.05%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "128.30"
$ws.Range("E48").Value = "  -0.62%  "

$ws.Range("E49").Value = "  +1.84%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.173"
$ws.Range("E50").Value = "  -0.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07101"
